# "Dons de la Lune finis" - complete the "Nouvelle Lune" gift rows (34-36)
# on the "Dons d'origine" sheet with the three missing Facettes: Chasseur
# Implacable, Diviser et Conquérir, and Brèche.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dons d'origine")

# --- Row 34 : Chasseur Implacable (•••) ---
$ws.Range("C34").Value = "Chasseur Implacable"
$ws.Range("E34").Value = "Peut importe à quel point la proie pense être en sécurité, l'Irraka atteint toujours sa cible.`nCette facette ne peut être activée que quand l'Irraka acquiers l'Etat Siskur-Dah."
$ws.Range("F34").Value = "E"
$ws.Range("G34").Value = "n/a"
$ws.Range("H34").Value = "Instantanée"
$ws.Range("I34").Value = "Toute la durée de la Siskur-Dah"
$ws.Range("J34").Value = "Pour toute la durée de la Facette, l'Irraka ajoute sa Renomée Ruse à ses jets suivants tant qu'elle le fait pour poursuivre la chasse.`n- Furtivité contre la proie.`n- Tentatives d'outrepasser n'importe quel obstacle ou sécurité la proie peut avoir mis entre elle et le danger, tels des jets de Larcin pour rentrer dans un immeuble où elle se cache`n- Des jets pour trouver des points d'accès aux trous dans lesquels la proie se tapis, des chemins alternatifs pour l'atteindre et autres moyens d'arriver à la proie malgré tous ses efforts"
$ws.Range("K34").Value = "n/a"
$ws.Range("L34").Value = "n/a"
$ws.Range("M34").Value = "n/a"
$ws.Range("N34").Value = "n/a"
$ws.Rows.Item(34).RowHeight = 180

# --- Row 35 : Diviser et Conquérir (••••) ---
$ws.Range("C35").Value = "Diviser et Conquérir"
$ws.Range("E35").Value = "L'Irraka n'a nul besoin d'attaquer les proies quand elles sont fortes et nombreuses. C'est bien mieux de les séparer avec des distractions et des appâts pour ensuite s'en occuper un par un."
$ws.Range("F35").Value = "E"
$ws.Range("G35").Value = "Manipulation + Subterfuge + Ruse contre Calme + Instinct Primal"
$ws.Range("H35").Value = "Contestée"
$ws.Range("I35").Value = "n/a"
$ws.Range("J35").Value = "L'Irraka peut utiliser cette Facette sur un seul individu qu'il peut voir et qui fait partie d'un groupe. L'Irraka doit fournir une distraction quelconque que ce soit un mouvement ou un son, mais elle peut être incroyablement mineure et seule la proie s'en rendra compte, initialement."
$ws.Range("K35").Value = 'La proie prends peur, elle gagne l''Etat "Spooked" et cherche desespérement à ne pas rester seule.'
$ws.Range("L35").Value = "La Facette échoue"
$ws.Range("M35").Value = 'La proie gagne l''Etat "Lured". Elle va volontairement se séparer du groupe et quitter leur présence immédiate pour investiguer, malgré son bon sens.'
$ws.Range("N35").Value = 'La Facette affecte aussi un nombre de ses compagnons égal à la Ruse de l''Irraka, leur donnant l''Etat "Lured" et les éparpillant à la suite de diversions imaginaires.'
$ws.Rows.Item(35).RowHeight = 90

# --- Row 36 : Brèche (•••••) ---
$ws.Range("C36").Value = "Brèche"
$ws.Range("E36").Value = "L'Irraka peut outrepasser à travers le Gantelet tel une ombre entre les deux mondes"
$ws.Range("F36").Value = "EEE"
$ws.Range("G36").Value = "Astuce + Furtivité + Ruse"
$ws.Range("H36").Value = "Instantanée"
$ws.Range("I36").Value = "n/a"
$ws.Range("J36").Value = "n/a"
$ws.Range("K36").Value = 'La tentative de l''Irraka d''ouvrir une brèche dans le Gantelet est violement repoussée. Elle souffre un Etat au choix parmis "Arm Wrack", "Leg Wrack" et "Stunned"'
$ws.Range("L36").Value = "La Facette échoue"
$ws.Range("M36").Value = "L'Irraka ouvre une brèche dans le Gantelet et atteint l'autre côté, arrivant à l'endroit correspondant en Chair ou Ombre. L'Irraka n'a pas besoin d'être à un locus pour ce faire."
$ws.Range("N36").Value = "La prochaine fois que l'Irraka utilise Brèche cette scène cela ne lui coûtera pas d'Essence."
$ws.Rows.Item(36).RowHeight = 60

# --- Update frozen-pane / active-cell selection to match the new content ---
$ws.Application.Goto($ws.Range("J37"), $true)
$ws.Range("J37").Select()
